$d = $word.ActiveDocument
$lb = [char]11

# --- Paragraph 1 (heading) & Paragraph 2 (empty spacer): add line spacing ---
$p1 = $d.Paragraphs(1)
$p1.Range.ParagraphFormat.LineSpacingRule = 1
$p1.Range.ParagraphFormat.LineSpacing = 18

$p2 = $d.Paragraphs(2)
$p2.Range.ParagraphFormat.LineSpacingRule = 1
$p2.Range.ParagraphFormat.LineSpacing = 18

# --- Paragraph 3: rewrite the whole intro paragraph ---
$p3 = $d.Paragraphs(3)
$s3 = $p3.Range.Start
$e3 = $p3.Range.End
$r3 = $d.Range($s3, $e3)
$r3.Text = "Hvert eneste år rejser en stor del af den danske befolkning på ferie i udlandet. Når man så sidder dér i flyveren på vej hjem, og længtes efter at komme hjem i vante omgivelser, er det nok de færreste som har fantasi til at forestille sig, at deres hjem i mellemtiden er blevet gennemrodet af indbrudstyve. " + $lb + "Ikke desto mindre er dette virkeligheden for rigtigt mange Danskere. " + $lb + "Man kommer hjem til boligen, og ser straks at døren til hjemmet er blevet brudt op, og får en meget ubehagelig følelse i kroppen. Efter hånden som man kommer igennem boligen, kan man se at alle værdier er blevet stjålet, TV og HIFI udstyr er væk, alle designer møbler er væk.  Det værste er dog dér hvor også uerstattelige genstande som billeder fra børnenes første år, eller den fantastiske hjemme-strikkede uro som lille Peter lavede i børnehaven, er blevet ødelagt. "

# --- Paragraph 4: statistics / project description paragraph ---
# Text before the first footnote reference.
$d.Content.Find.Execute("anmeldt ikke mindre end 32.974", $false, $false, $false, $false, $false, $true, 1, $false, "anmeldt 32.974", 2) | Out-Null

# Text between footnote 1 and footnote 2.
$d.Content.Find.Execute("indbrud i beboelses ejendomme rundt omkring i Danmark", $false, $false, $false, $false, $false, $true, 1, $false, "indbrud i beboelses ejendomme i Danmark", 2) | Out-Null

# Text from right after footnote 2 up to (but excluding) the "home" word (keep the
# proofErr-wrapped "home" run untouched).
$f = $d.Content
$f.Find.Execute(" røverier. ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$segStart = $f.End
$g = $d.Content
$g.Find.Execute("home", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$segEnd = $g.Start
$seg = $d.Range($segStart, $segEnd)
$seg.Text = $lb + "Det er altså tydeligt, at man kan mindske risikoen for at komme hjem til en sådan ubehagelig overraskelse, hvis man kan få indbrudstyvene til at tro at der er nogen hjemme i boligen, selvom man er afsted på ferie." + $lb + "Vi vil i vores projekt skabe et system som kan netop dette. " + $lb + "Baseret på ”"

# Text right after "home" up to the end of the "kodelås." sentence.
$h = $d.Content
$h.Find.Execute("home", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$segStart2 = $h.End
$i2 = $d.Content
$i2.Find.Execute("tidsplan, som gemmes i systemets hukommelse. Tidsplanen skal beskyttes af en kodelås, som skal indtastes inden man via. den grafiske brugerflade på en tilkoblet PC, kan lave ændringer i tidsplanen, for at forhindre at hvem som helst kan lave ændringer i de gemte indstillinger.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$segEnd2 = $i2.End
$seg2 = $d.Range($segStart2, $segEnd2)
$seg2.Text = "-automation” konceptet, vil vi bygge et system som kan få det til at se ud som om der er folk hjemme, ved at tænde og slukke for lyset i de forskellige rum i hjemmet, eller for en enkelt lampe. " + $lb + "Systemet vil benytte det eksisterende lysnet i boligen til kommunikationen, så man som bruger blot kan tilslutte systemet til stikkontakten, og så virker det.  Systemet skal kunne kører automatisk ud fra en bruger styret tidsplan. Tidsplanen skal kunne ændres gennem en grafisk brugerflade via. en tilkoblet PC, efter korrekt indtastning af kode på systemets kodelås."

# Remaining tail of the paragraph (runs through the _GoBack bookmark, which the engine
# keeps alive and re-anchors at the end of the replaced range).
$j = $d.Content
$j.Find.Execute("Systemet vil benytte det indlagte lysnet i huset til kommunikationen", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$segStart3 = $j.Start
$p4 = $d.Paragraphs(4)
$segEnd3 = $p4.Range.End - 1
$seg3 = $d.Range($segStart3, $segEnd3)
$seg3.Text = $lb + "Systemet skal være i stand til at registrere hvis der opstår kommunikationsfejl mellem styrerbokse og  enheder og logge disse. Brugeren skal så få en oversigt over fejl der har været siden sidste login, når man igen tilkobler og logger ind via. PC. "

# Add a trailing line break after the bookmark, at the end of the paragraph.
$bm2 = $d.Bookmarks("_GoBack")
$p4b = $d.Paragraphs(4)
$tail = $d.Range($bm2.Range.End, $p4b.Range.End - 1)
$tail.InsertAfter($lb)

# --- Footnotes: fix "statestik" -> "statistik" typo and drop the now-unneeded proofErr marks ---
$fn1 = $d.Footnotes(1)
$fn1.Range.Text = " Tal fra Danmarks statistik: http://www.dst.dk/da/Statistik/NytHtml?cid=20617"
$fn2 = $d.Footnotes(2)
$fn2.Range.Text = " Tal fra Danmarks statistik: http://www.dst.dk/da/Statistik/NytHtml?cid=20617"

Write-Output "done"
